$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("B3").Value = "MTRM"
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "Circuitos Elétricos 2"
$ws.Range("E3").Value = "-"
$ws.Range("F3").Value = "CAD"

# Row 4
$ws.Range("B4").Value = "MTRM"
$ws.Range("C4").Value = "Acionamentos"
$ws.Range("D4").Value = "Circuitos Elétricos 2"
$ws.Range("E4").Value = "Programação"
$ws.Range("F4").Value = "EAP"

# Row 6
$ws.Range("B6").Value = "Sistemas digitais"
$ws.Range("C6").Value = "EAP"
$ws.Range("D6").Value = "Sistemas digitais"
